$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(35, "nico@gmail.com", "Eliminó producto: Yuval (ID 3)", "2026-02-16 17:02:51"),
    @(36, "nico@gmail.com", "Inicio de sesión exitoso", "2026-02-16 17:16:04"),
    @(37, "nico@gmail.com", "Inicio de sesión exitoso", "2026-02-16 22:34:53"),
    @(38, "nico@gmail.com", "Inicio de sesión exitoso", "2026-02-16 23:25:42"),
    @(39, "nico@gmail.com", "Inicio de sesión exitoso", "2026-02-17 00:32:37"),
    @(40, "nico@gmail.com", "Inicio de sesión exitoso", "2026-02-17 00:42:13"),
    @(41, "nico@gmail.com", "Inicio de sesión exitoso", "2026-02-17 13:03:31"),
    @(42, "nico@gmail.com", "Inicio de sesión exitoso", "2026-02-17 13:34:38"),
    @(43, "nico@gmail.com", "Inicio de sesión exitoso", "2026-02-17 14:37:50"),
    @(44, "nico@gmail.com", "Inicio de sesión exitoso", "2026-02-17 14:38:12")
)

$startRow = 36
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
}
